$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (IntervarClassification), shifting
# columns E:O to F:P. xlShiftToRight = -4161, xlFormatFromRightOrBelow = 0
# so the new column inherits formatting from the (old) column E, now F.
$ws.Range("E1:E3").EntireColumn.Insert(-4161, 0)

# Populate the new column E with header + data.
$ws.Range("E1").Value = "IntervarConsequence"
$ws.Range("E2").Value = "exonic,frameshift deletion"
$ws.Range("E3").Value = "exonic,frameshift deletion"
